$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the previously used range entirely before rewriting the documented
# workflow-tree content (the activity tree was renumbered/restructured).
$ws.Range("A1:V118").Clear()

$ws.Range("A1").Value = "1.1 Extract_Fields (Sequence)"
$ws.Range("B2").Value = "Private = False"
$ws.Range("B3").Value = "Activities"
$ws.Range("C4").Value = "1.2 For Each (ForEach<String>)"
$ws.Range("D5").Value = "Values = Directory.GetFiles(`"Invoices`",`"*.pdf`")"
$ws.Range("D6").Value = "Private = False"
$ws.Range("D7").Value = "TypeArgument = System.String"
$ws.Range("D8").Value = "Body"
$ws.Range("E9").Value = "1.3 Body (Sequence)"
$ws.Range("F10").Value = "Private = False"
$ws.Range("F11").Value = "Activities"
$ws.Range("G12").Value = "1.75 Start Process - Adobe Acrobat (StartProcess)"
$ws.Range("H13").Value = "FileName = invoice"
$ws.Range("H14").Value = "Private = False"
$ws.Range("G15").Value = "1.4 Attach Window - Invoice.pdf (WindowScope)"
$ws.Range("H16").Value = "Selector = <wnd app='acrobat.exe' cls='AcrobatSDIWindow' title='*.pdf - Adobe Acrobat Standard DC' />"
$ws.Range("H17").Value = "Private = False"
$ws.Range("H18").Value = "Body"
$ws.Range("I19").Value = "1.5 Do (Sequence)"
$ws.Range("J20").Value = "Private = False"
$ws.Range("J21").Value = "Variables"
$ws.Range("K22").Value = "textInvoice(GenericValue)"
$ws.Range("K23").Value = "textName(GenericValue)"
$ws.Range("K24").Value = "textDate(GenericValue)"
$ws.Range("K25").Value = "dataTable(DataTable)"
$ws.Range("K26").Value = "counter(Int32)"
$ws.Range("J27").Value = "Activities"
$ws.Range("K28").Value = "1.73 Maximize Window - Adobe Acrobat (MaximizeWindow)"
$ws.Range("L29").Value = "Private = False"
$ws.Range("K30").Value = "1.68 Get Text - Customer Name (GetValue)"
$ws.Range("L31").Value = "Value = textName"
$ws.Range("L32").Value = "Target"
$ws.Range("M33").Value = "Selector = <wnd app='acrobat.exe' cls='AcrobatSDIWindow' title='*.pdf - Adobe Acrobat Standard DC' /><wnd aaname='Document Pane' cls='AVL_AVView' title='AVScrolledPageView' /><wnd cls='AVL_AVView' title='AVPageView' /><ctrl name='*, ' role='text' />"
$ws.Range("L34").Value = "Private = False"
$ws.Range("K35").Value = "1.63 Assign - Remove , (Assign)"
$ws.Range("L36").Value = "To = textName"
$ws.Range("L37").Value = "Value = textName.Replace(`",`",`" `")"
$ws.Range("L38").Value = "Private = False"
$ws.Range("K39").Value = "1.58 Assign - Trim (Assign)"
$ws.Range("L40").Value = "To = textName"
$ws.Range("L41").Value = "Value = textName.Trim"
$ws.Range("L42").Value = "Private = False"
$ws.Range("K43").Value = "1.53 Get Text - Invoice # (GetValue)"
$ws.Range("L44").Value = "Value = textInvoice"
$ws.Range("L45").Value = "Target"
$ws.Range("M46").Value = "Selector = <wnd app='acrobat.exe' cls='AcrobatSDIWindow' title='*.pdf - Adobe Acrobat Standard DC' /><wnd aaname='Document Pane' cls='AVL_AVView' title='AVScrolledPageView' /><wnd cls='AVL_AVView' title='AVPageView' /><ctrl idx='2' role='row' /><ctrl name='Invoice # * ' role='text' />"
$ws.Range("L47").Value = "Private = False"
$ws.Range("K48").Value = "1.48 Get Text - Invoice Date (GetValue)"
$ws.Range("L49").Value = "Value = textDate"
$ws.Range("L50").Value = "Target"
$ws.Range("M51").Value = "Selector = <wnd app='acrobat.exe' cls='AcrobatSDIWindow' title='*.pdf - Adobe Acrobat Standard DC' /><wnd aaname='Document Pane' cls='AVL_AVView' title='AVScrolledPageView' /><wnd cls='AVL_AVView' title='AVPageView' /><ctrl idx='2' role='row' /><ctrl name='Invoice Date: * ' role='text' />"
$ws.Range("L52").Value = "Private = False"
$ws.Range("K53").Value = "1.9 Do While - Grab Table (DoWhile)"
$ws.Range("L54").Value = "Condition = counter <4"
$ws.Range("L55").Value = "Private = False"
$ws.Range("L56").Value = "Body"
$ws.Range("M57").Value = "1.10 Sequence - Grab Table (Sequence)"
$ws.Range("N58").Value = "Private = False"
$ws.Range("N59").Value = "Activities"
$ws.Range("O60").Value = "1.38 Extract Structured Data - Order Information (ExtractData)"
$ws.Range("P61").Value = "ExtractMetadata = <extract-table get_columns_name='1' get_empty_columns='1' columns_name_source='Longest' />"
$ws.Range("P62").Value = "MaxNumberOfResults = 100"
$ws.Range("P63").Value = "DataTable = dataTable"
$ws.Range("P64").Value = "SimulateClick = True"
$ws.Range("P65").Value = "Target"
$ws.Range("Q66").Value = "Selector = <wnd aaname='Document Pane' cls='AVL_AVView' title='AVScrolledPageView' /><wnd cls='AVL_AVView' title='AVPageView' /><ctrl idx='{{counter}}' role='table' />"
$ws.Range("P67").Value = "ContinueOnError = True"
$ws.Range("P68").Value = "Private = False"
$ws.Range("O69").Value = "1.16 If - counter = 2 (If)"
$ws.Range("P70").Value = "Condition = counter=2"
$ws.Range("P71").Value = "Private = False"
$ws.Range("P72").Value = "Then"
$ws.Range("Q73").Value = "1.28 Excel Application Scope - Create or Open Workbook (ExcelApplicationScope)"
$ws.Range("R74").Value = "WorkbookPath = `"Excel Invoices\`"+textName+`"_`"+textInvoice+`".xlsx`""
$ws.Range("R75").Value = "Visible = True"
$ws.Range("R76").Value = "CreateNewFile = True"
$ws.Range("R77").Value = "AutoSave = True"
$ws.Range("R78").Value = "ReadOnly = False"
$ws.Range("R79").Value = "MacroSetting = EnableAll"
$ws.Range("R80").Value = "Private = False"
$ws.Range("R81").Value = "Body"
$ws.Range("S82").Value = "1.29 Do (Sequence)"
$ws.Range("T83").Value = "Private = False"
$ws.Range("T84").Value = "Activities"
$ws.Range("U85").Value = "1.30 Write Range - Write Data Table (ExcelWriteRange)"
$ws.Range("V86").Value = "StartingCell = A3"
$ws.Range("V87").Value = "DataTable = dataTable"
$ws.Range("V88").Value = "AddHeaders = True"
$ws.Range("V89").Value = "SheetName = Sheet1"
$ws.Range("V90").Value = "Private = False"
$ws.Range("P91").Value = "Else"
$ws.Range("Q92").Value = "1.19 Excel Application Scope - Open Workbook (ExcelApplicationScope)"
$ws.Range("R93").Value = "WorkbookPath = `"Excel Invoices\`"+textName+`"_`"+textInvoice+`".xlsx`""
$ws.Range("R94").Value = "Visible = True"
$ws.Range("R95").Value = "CreateNewFile = True"
$ws.Range("R96").Value = "AutoSave = True"
$ws.Range("R97").Value = "ReadOnly = False"
$ws.Range("R98").Value = "MacroSetting = EnableAll"
$ws.Range("R99").Value = "Private = False"
$ws.Range("R100").Value = "Body"
$ws.Range("S101").Value = "1.20 Do (Sequence)"
$ws.Range("T102").Value = "Private = False"
$ws.Range("T103").Value = "Activities"
$ws.Range("U104").Value = "1.21 Append Range - Write Next Page (ExcelAppendRange)"
$ws.Range("V105").Value = "DataTable = dataTable"
$ws.Range("V106").Value = "SheetName = Sheet1"
$ws.Range("V107").Value = "Private = False"
$ws.Range("O108").Value = "1.11 Assign - counter+1 (Assign)"
$ws.Range("P109").Value = "To = counter"
$ws.Range("P110").Value = "Value = counter+1"
$ws.Range("P111").Value = "Private = False"
$ws.Range("K112").Value = "1.8 Close Application - Adobe Acrobat (CloseApplication)"
$ws.Range("L113").Value = "Target"
$ws.Range("L114").Value = "Private = False"
